$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlinks (EMAIL column used to be E, with mailto: links) ---
$ws.Hyperlinks.Delete()

# --- Header row ---
$ws.Range("A1").Value = "Student ID"
$ws.Range("C1").Value = "EMAIL"
$ws.Range("D1").Value = "GENDER"
$ws.Range("E1").Value = "MAJOR"

# --- Row 2: Diep Truong Khanh Bang ---
$ws.Range("A2").Value = "'52200238"
$ws.Range("B2").Value = "Diep Truong Khanh Bang"
$ws.Range("C2").Formula = "=A2&""@student.tdtu.edu.vn"""
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("D2").Value = "Female"
$ws.Range("E2").Value = "Mang may tinh va truyen thong du lieu"

# --- Row 3: Tang Duy Hao ---
$ws.Range("A3").Value = "'52200210"
$ws.Range("B3").Value = "Tang Duy Hao"
$ws.Range("C3").Formula = "=A3&""@student.tdtu.edu.vn"""
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("D3").Value = "Male"
$ws.Range("E3").Value = "Khoa hoc may tinh"

# --- Row 4: Ho Bao Ngan (new row) ---
$ws.Range("A4").Value = "'52200243"
$ws.Range("B4").Value = "Ho Bao Ngan"
$ws.Range("C4").Formula = "=A4&""@student.tdtu.edu.vn"""
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("D4").Value = "Female"
$ws.Range("E4").Value = "Ky thuat phan mem"

# --- Drop the now-unused USERNAME/PASSWORD .. ROLE columns (old C,D,F,G) ---
# After placing the new A:E layout, the leftover data lives in columns F and G.
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(6).Delete()

# --- Cosmetics: font, column widths, row heights ---
$ws.Range("A1").Font.Name = "Arial"
$wb.Styles("Normal").Font.Name = "Arial"
$wb.Styles("Hyperlink").Font.Name = "Arial"

$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 20.833333333333332
$ws.Columns.Item(3).ColumnWidth = 25.5
$ws.Columns.Item(5).ColumnWidth = 31

$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8
$ws.Rows.Item(8).RowHeight = 13.8

# --- Selection like the saved file ---
$ws.Range("E8").Select()
